# edit.ps1
# Updates the cryptos worksheet with refreshed coin ranking data.
# Corresponds to commit: "Updated cryptos list on Fri Jun 16 23:48:43 UTC 2023 with GitHub Actions"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @('Bitcoin', 'https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc', '26.340.60', '  +3.09%  '),
    @('Ethereum', 'https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth', '1.716.88', '  +3.18%  '),
    @('TetherUSD', 'https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt', '0.9995', '  +0.03%  '),
    @('BNB', 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb', '239.16', '  +1.40%  '),
    @('USDC', 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc', '1.000', '  +0.02%  '),
    @('XRP', 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp', '0.4751', '  -0.90%  '),
    @('OKB', 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb', '41.13', '  +3.12%  '),
    @('Cardano', 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada', '0.2632', '  +0.72%  '),
    @('Dogecoin', 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge', '0.06213', '  +1.06%  '),
    @('WrappedEther', 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth', '1.714.48', '  +3.07%  '),
    @('TRON', 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx', '0.07055', '  -0.49%  '),
    @('Solana', 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol', '15.32', '  +4.00%  '),
    @('Polkadot', 'https://coinranking.com/coin/25W7FG7om+polkadot-dot', '4.419', '  +1.21%  '),
    @('Polygon', 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic', '0.5890', '  -0.28%  '),
    @('Litecoin', 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc', '76.02', '  +2.21%  '),
    @('Dai', 'https://coinranking.com/coin/MoTuySvg7+dai-dai', '1.000', '  +0.00%  '),
    @('BinanceUSD', 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd', '1.000', '  +0.06%  '),
    @('WrappedBTC', 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc', '26.334.70', '  +3.12%  '),
    @('ShibaInu', 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib', '0.000006816', '  +1.00%  '),
    @('Avalanche', 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax', '11.55', '  +1.38%  '),
    @('WrappedliquidstakedEther2.0', 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth', '1.935.85', '  +2.99%  '),
    @('Uniswap', 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni', '4.545', '  +2.76%  '),
    @('Cosmos', 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom', '8.756', '  +1.19%  '),
    @('Chainlink', 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link', '5.321', '  +0.41%  '),
    @('Monero', 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr', '134.15', '  -0.18%  '),
    @('EthereumClassic', 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc', '15.24', '  +1.40%  '),
    @('Toncoin', 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton', '1.402', '  +0.37%  '),
    @('BitcoinCash', 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch', '108.00', '  +3.23%  '),
    @('LidoDAOToken', 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo', '1.754', '  +4.18%  '),
    @('InternetComputer(DFINITY)', 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp', '3.997', '  +1.06%  '),
    @('Filecoin', 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil', '3.690', '  +1.11%  '),
    @('Stellar', 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm', '0.07746', '  +1.68%  '),
    @('Hedera', 'https://coinranking.com/coin/jad286TjB+hedera-hbar', '0.04441', '  +3.06%  '),
    @('HuobiToken', 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht', '2.614', '  -0.15%  '),
    @('ARBITRUM', 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb', '0.9754', '  +2.90%  '),
    @('ImmutableX', 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx', '0.6188', '  +1.11%  '),
    @('TrustWalletToken', 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt', '0.9260', '  +9.20%  '),
    @('Quant', 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt', '112.16', '  +14.69%  '),
    @('MXToken', 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx', '2.421', '  -7.17%  '),
    @('RenderToken', 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr', '1.919', '  +2.49%  '),
    @('PaxDollar', 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp', '1.000', '  +0.01%  '),
    @('VeChain', 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet', '0.01470', '  -1.61%  '),
    @('FraxShare', 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs', '5.323', '  +13.48%  '),
    @('TheSandbox', 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand', '0.3810', '  +1.44%  '),
    @('Algorand', 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo', '0.1162', '  +3.84%  '),
    @('Aptos', 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt', '6.297', '  +1.44%  '),
    @('Cronos', 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro', '0.05282', '  +0.43%  '),
    @('Elrond', 'https://coinranking.com/coin/omwkOTglq+elrond-egld', '30.33', '  +3.32%  '),
    @('EnergySwap', 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens', '7.687', '  +5.26%  '),
    @('Aave', 'https://coinranking.com/coin/ixgUfzmLR+aave-aave', '50.72', '  +1.17%  ')
)

# Force column D (Price) to be treated as text so that values such as
# "1.000", "0.5890" or "26.334.70" keep their exact original formatting
# instead of being auto-converted to numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $rec = $data[$i]
    $ws.Cells.Item($row, 2).Value = $rec[0]
    $ws.Cells.Item($row, 3).Value = $rec[1]
    $ws.Cells.Item($row, 4).Value = $rec[2]
    $ws.Cells.Item($row, 5).Value = $rec[3]
}
